$d = $word.ActiveDocument

# --- Edit 1 -------------------------------------------------------------
# "... odděleny znakem = (rovná se)." -> "... odděleny znakem _ (podtržítko)."
# The tail gets split into several separate runs, mirroring the diff.
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute("= (rovná se).", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Edit 1: target text not found"
}

$tailStart = $rng1.Start
$tailEnd = $rng1.End

# Remove the old tail text entirely, leaving an insertion point behind.
$tailRange = $d.Range($tailStart, $tailEnd)
$tailRange.Text = ""

# Re-insert the replacement as five distinct runs: "_", " (", "podtržítko", ")."
$p1 = $d.Range($tailStart, $tailStart)
$p1.InsertAfter("_")

$p2 = $d.Range($p1.End, $p1.End)
$p2.InsertAfter(" (")

$p3 = $d.Range($p2.End, $p2.End)
$p3.InsertAfter("podtržítko")

$p4 = $d.Range($p3.End, $p3.End)
$p4.InsertAfter(").")

# --- Edit 2 -------------------------------------------------------------
# "Při spuštění nové hry se nastaví ..." -> "nové hry" gets the Nadpis5Char
# character style, splitting the run into three.
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute("nové hry", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Edit 2: target text not found"
}
$rng2.Style = "Nadpis5Char"
